$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Addr, [string]$Text)
    $cell = $ws.Range($Addr)
    $needsQuote = $Text -match '^[+-]?[0-9]*\.?[0-9]+$'
    if ($needsQuote) {
        $cell.Value = "'" + $Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Text
    }
}

Set-CellText "D2" "48.204.96"
Set-CellText "E2" "  +0.16%  "
Set-CellText "D3" "2.501.95"
Set-CellText "E3" "  -0.92%  "
Set-CellText "E4" "  -0.03%  "
Set-CellText "D5" "318.13"
Set-CellText "E5" "  -1.77%  "
Set-CellText "D6" "106.07"
Set-CellText "E6" "  -2.69%  "
Set-CellText "E7" "  -1.60%  "
Set-CellText "E8" "  -0.02%  "
Set-CellText "E9" "  -3.80%  "
Set-CellText "D10" "38.95"
Set-CellText "E10" "  -4.07%  "
Set-CellText "D11" "20.26"
Set-CellText "E11" "  -0.69%  "
Set-CellText "E12" "  -2.70%  "
Set-CellText "E13" "  -0.02%  "
Set-CellText "E14" "  -2.33%  "
Set-CellText "D15" "2.893.24"
Set-CellText "E15" "  -0.88%  "
Set-CellText "D16" "2.504.16"
Set-CellText "E16" "  -0.84%  "
Set-CellText "E17" "  -3.34%  "
Set-CellText "D18" "48.065.68"
Set-CellText "E18" "  +0.19%  "
Set-CellText "D19" "3.00"
Set-CellText "E19" "  +10.82%  "
Set-CellText "E20" "  -2.90%  "
Set-CellText "D21" "6.60"
Set-CellText "E21" "  -0.52%  "
Set-CellText "E22" "  -2.31%  "
Set-CellText "E23" "  -1.59%  "
Set-CellText "D24" "268.55"
Set-CellText "E24" "  -0.48%  "
Set-CellText "D25" "2.52"
Set-CellText "E25" "  -2.38%  "
Set-CellText "E26" "  +0.20%  "
Set-CellText "D27" "25.83"
Set-CellText "E27" "  -1.60%  "
Set-CellText "D28" "2.28"
Set-CellText "E28" "  +8.26%  "
Set-CellText "E29" "  -4.13%  "
Set-CellText "D30" "0.139"
Set-CellText "E30" "  -4.61%  "
Set-CellText "D31" "34.61"
Set-CellText "E31" "  -3.14%  "
Set-CellText "D32" "49.36"
Set-CellText "E32" "  -0.92%  "
Set-CellText "B33" "Celestia"
Set-CellText "C33" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-CellText "D33" "19.19"
Set-CellText "E33" "  -3.93%  "
Set-CellText "B34" "FirstDigitalUSD"
Set-CellText "C34" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText "D34" "1.00"
Set-CellText "E34" "  +0.01%  "
Set-CellText "E35" "  -2.08%  "
Set-CellText "D36" "0.0774"
Set-CellText "E36" "  -2.63%  "
Set-CellText "E37" "  -2.45%  "
Set-CellText "E38" "  -3.06%  "
Set-CellText "D39" "2.89"
Set-CellText "E39" "  -3.66%  "
Set-CellText "D40" "122.87"
Set-CellText "E40" "  +3.22%  "
Set-CellText "B41" "EnergySwap"
Set-CellText "C41" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D41" "22.53"
Set-CellText "E41" "  +0.34%  "
Set-CellText "B42" "Stellar"
Set-CellText "C42" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D42" "0.111"
Set-CellText "E42" "  -1.44%  "
Set-CellText "E43" "  +0.99%  "
Set-CellText "D44" "0.0302"
Set-CellText "E44" "  +0.68%  "
Set-CellText "D45" "2.002.51"
Set-CellText "E45" "  -0.47%  "
Set-CellText "E46" "  +0.47%  "
Set-CellText "E47" "  +1.52%  "
Set-CellText "E48" "  -2.78%  "
Set-CellText "D49" "8.94"
Set-CellText "E49" "  -2.24%  "
Set-CellText "D50" "5.22"
Set-CellText "E50" "  -0.72%  "
Set-CellText "D51" "78.93"
Set-CellText "E51" "  -1.09%  "
